# GUI.py change: rename the login-details table from Employee_ID/Username/
# Password to a simpler User_ID/Password pair, and refresh the sample
# credentials (drop the old Username column entirely).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "User_ID"
$ws.Range("B1").Value = "Password"

# Data row - A2 ("978232") looks numeric, so force text entry (matching the
# original file, where this row was stored as shared-string text, not a
# number) and then restore the default "Normal" style so no stray
# number-format style sticks to the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "978232"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "ijdf9f"

# The old third column (Password/3rqw) is removed entirely.
$ws.Range("C1:C2").Delete()
